$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.774.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.534.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.484"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.24"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0580"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.751.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.535.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.506"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.752.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0680"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.60%  "
$ws.Range("E24").Value = "  -3.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0454"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.361.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.958"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.73%  "
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.519"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.22%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.801"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.665.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0504"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0976"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0942"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.14%  "
